# Last Song formatting fix.
#
# Previously, the closing song ("Palabras de conclusión(3 mins.)|Canción Ny
# oración") was glued onto whatever cell happened to hold the last numbered
# point of the "NUESTRA VIDA CRISTIANA" block - so the song ended up on a
# different row depending on how many points a given week had.
#
# Now the song is parsed out of that text and written on its own, as a plain
# "Canción N" cell, one row below the last point - i.e. at the true end of
# the NVC block, growing the sheet by a row when needed. Weeks that don't
# follow the usual pattern (e.g. a "weird" Conmemoración week) are left
# untouched instead of being guessed at.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ur = $ws.UsedRange
$firstRow = $ur.Row()
$firstCol = $ur.Column()
$lastRow = $firstRow + $ur.Rows.Count() - 1
$lastCol = $firstCol + $ur.Columns.Count() - 1

# Section header that the closing song always sits somewhere below, in the
# same (weekly) column.
$nvcLabel = "NUESTRA VIDA CRISTIANA"

for ($col = $firstCol; $col -le $lastCol; $col++) {

    # Locate the NVC header in this column.
    $nvcRow = -1
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        if ($ws.Cells.Item($r, $col).Value() -eq $nvcLabel) {
            $nvcRow = $r
        }
    }

    if ($nvcRow -eq -1) {
        # No NVC block in this column at all - nothing to fix.
        continue
    }

    # Scan down from the header for the old glued-together closing line,
    # recognisable by the "Palabras de conclusión ... | Canción ..." pipe.
    $songRow = -1
    $songText = ""
    for ($r = $nvcRow + 1; $r -le $lastRow; $r++) {
        $cellVal = $ws.Cells.Item($r, $col).Value()
        if ($cellVal -ne $null -and $cellVal.Contains("|")) {
            $songRow = $r
            $songText = $cellVal
        }
    }

    if ($songRow -eq -1) {
        # "Weird" week (e.g. Conmemoración) that doesn't match the usual
        # closing-song pattern - skip it rather than mangle it.
        continue
    }

    # Pull just the song number out of "...|Canción Ny oración".
    if ($songText -notmatch 'Canción\s*(\d+)') {
        continue
    }
    $songNumber = $matches[1]

    # Remove the old combined text and place the bare song by itself one
    # row down, at the real end of the block.
    $ws.Cells.Item($songRow, $col).Value = ""
    $ws.Cells.Item($songRow + 1, $col).Value = "Canción $songNumber"
}
